# Auto-generated: apply row permutation for rows 39-49 (record reshuffle)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39 now holds the record previously at row 48
$ws.Range("A39").Value = 111880484
$ws.Range("B39").Value = 90658
$ws.Range("E39").Value = 4361
$ws.Range("F39").Value = 'Orange taggsvamp'
$ws.Range("G39").Value = 'Hydnellum aurantiacum'
$ws.Range("H39").Value = '(Batsch:Fr.) P.Karst.'
$ws.Range("I39").NumberFormat = "@"
$ws.Range("I39").Value = "11"
$ws.Range("Q39").Value = 509900.7891887496
$ws.Range("R39").Value = 6753525.142772059
$ws.Range("AJ39").Value = 'tall'
$ws.Range("AK39").Value = 'Pinus sylvestris'
$ws.Range("AO39").Value = 'Pinus sylvestris'

# Row 40 now holds the record previously at row 49
$ws.Range("A40").Value = 111880601
$ws.Range("B40").Value = 88966
$ws.Range("E40").Value = 5754
$ws.Range("F40").Value = 'Gultoppig fingersvamp'
$ws.Range("G40").Value = 'Ramaria testaceoflava'
$ws.Range("H40").Value = '(Bres.) Corner'
$ws.Range("I40").NumberFormat = "@"
$ws.Range("I40").Value = "4"
$ws.Range("Q40").Value = 509941.5744066621
$ws.Range("R40").Value = 6753224.672924293
$ws.Range("AJ40").Value = 'tall'
$ws.Range("AK40").Value = 'Pinus sylvestris'
$ws.Range("AO40").Value = 'Pinus sylvestris'

# Row 41 now holds the record previously at row 42
$ws.Range("A41").Value = 111880509
$ws.Range("B41").Value = 90652
$ws.Range("E41").Value = 3100
$ws.Range("F41").Value = 'Talltaggsvamp'
$ws.Range("G41").Value = 'Bankera fuligineoalba'
$ws.Range("H41").Value = '(Schmidt : Fr.) Pouzar'
$ws.Range("I41").NumberFormat = "@"
$ws.Range("I41").Value = "6"
$ws.Range("Q41").Value = 509834.2096935506
$ws.Range("R41").Value = 6753644.114383955
$ws.Range("AJ41").Value = 'tall'
$ws.Range("AK41").Value = 'Pinus sylvestris'
$ws.Range("AO41").Value = 'Pinus sylvestris'

# Row 42 now holds the record previously at row 47
$ws.Range("A42").Value = 111880562
$ws.Range("B42").Value = 90658
$ws.Range("E42").Value = 4361
$ws.Range("F42").Value = 'Orange taggsvamp'
$ws.Range("G42").Value = 'Hydnellum aurantiacum'
$ws.Range("H42").Value = '(Batsch:Fr.) P.Karst.'
$ws.Range("I42").NumberFormat = "@"
$ws.Range("I42").Value = "3"
$ws.Range("Q42").Value = 509657.7198006394
$ws.Range("R42").Value = 6753521.069647122
$ws.Range("AJ42").Value = 'tall'
$ws.Range("AK42").Value = 'Pinus sylvestris'
$ws.Range("AO42").Value = 'Pinus sylvestris'

# Row 43 now holds the record previously at row 45
$ws.Range("A43").Value = 111880500
$ws.Range("B43").Value = 88966
$ws.Range("E43").Value = 5754
$ws.Range("F43").Value = 'Gultoppig fingersvamp'
$ws.Range("G43").Value = 'Ramaria testaceoflava'
$ws.Range("H43").Value = '(Bres.) Corner'
$ws.Range("I43").NumberFormat = "@"
$ws.Range("I43").Value = "4"
$ws.Range("Q43").Value = 509899.1991435916
$ws.Range("R43").Value = 6753571.34232254
$ws.Range("AJ43").Value = 'gran'
$ws.Range("AK43").Value = 'Picea abies'
$ws.Range("AO43").Value = 'Picea abies'

# Row 45 now holds the record previously at row 43
$ws.Range("A45").Value = 111880574
$ws.Range("B45").Value = 90658
$ws.Range("E45").Value = 4361
$ws.Range("F45").Value = 'Orange taggsvamp'
$ws.Range("G45").Value = 'Hydnellum aurantiacum'
$ws.Range("H45").Value = '(Batsch:Fr.) P.Karst.'
$ws.Range("I45").NumberFormat = "@"
$ws.Range("I45").Value = "2"
$ws.Range("Q45").Value = 509595.7160662179
$ws.Range("R45").Value = 6753391.52735021
$ws.Range("AJ45").Value = 'tall'
$ws.Range("AK45").Value = 'Pinus sylvestris'
$ws.Range("AO45").Value = 'Pinus sylvestris'

# Row 46 now holds the record previously at row 40
$ws.Range("A46").Value = 111880532
$ws.Range("B46").Value = 90652
$ws.Range("E46").Value = 3100
$ws.Range("F46").Value = 'Talltaggsvamp'
$ws.Range("G46").Value = 'Bankera fuligineoalba'
$ws.Range("H46").Value = '(Schmidt : Fr.) Pouzar'
$ws.Range("I46").NumberFormat = "@"
$ws.Range("I46").Value = "2"
$ws.Range("Q46").Value = 509682.5105515064
$ws.Range("R46").Value = 6753540.591470475
$ws.Range("AJ46").Value = 'tall'
$ws.Range("AK46").Value = 'Pinus sylvestris'
$ws.Range("AO46").Value = 'Pinus sylvestris'

# Row 47 now holds the record previously at row 46
$ws.Range("A47").Value = 111880591
$ws.Range("B47").Value = 90658
$ws.Range("E47").Value = 4361
$ws.Range("F47").Value = 'Orange taggsvamp'
$ws.Range("G47").Value = 'Hydnellum aurantiacum'
$ws.Range("H47").Value = '(Batsch:Fr.) P.Karst.'
$ws.Range("I47").NumberFormat = "@"
$ws.Range("I47").Value = "8"
$ws.Range("Q47").Value = 509822.1902239832
$ws.Range("R47").Value = 6753234.069152902
$ws.Range("AJ47").Value = 'tall'
$ws.Range("AK47").Value = 'Pinus sylvestris'
$ws.Range("AO47").Value = 'Pinus sylvestris'

# Row 48 now holds the record previously at row 41
$ws.Range("A48").Value = 111880475
$ws.Range("B48").Value = 88966
$ws.Range("E48").Value = 5754
$ws.Range("F48").Value = 'Gultoppig fingersvamp'
$ws.Range("G48").Value = 'Ramaria testaceoflava'
$ws.Range("H48").Value = '(Bres.) Corner'
$ws.Range("I48").NumberFormat = "@"
$ws.Range("I48").Value = "2"
$ws.Range("Q48").Value = 509957.7514087428
$ws.Range("R48").Value = 6753362.853637428
$ws.Range("AJ48").Value = 'gran'
$ws.Range("AK48").Value = 'Picea abies'
$ws.Range("AO48").Value = 'Picea abies'

# Row 49 now holds the record previously at row 39
$ws.Range("A49").Value = 111880580
$ws.Range("B49").Value = 90658
$ws.Range("E49").Value = 4361
$ws.Range("F49").Value = 'Orange taggsvamp'
$ws.Range("G49").Value = 'Hydnellum aurantiacum'
$ws.Range("H49").Value = '(Batsch:Fr.) P.Karst.'
$ws.Range("I49").NumberFormat = "@"
$ws.Range("I49").Value = "3"
$ws.Range("Q49").Value = 509755.441071702
$ws.Range("R49").Value = 6753236.317390828
$ws.Range("AJ49").Value = 'tall'
$ws.Range("AK49").Value = 'Pinus sylvestris'
$ws.Range("AO49").Value = 'Pinus sylvestris'

